$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "'289.57"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'-3.83%"
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).Value = "'30.94"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'-3.77%"
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,4).Value = "'4.867"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'-2.26%"
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).Value = "'0.07126"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'-9.63%"
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).Value = "'1.807"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'-14.11%"
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).Value = "'7.651"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'-1.81%"
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).Value = "'3.780"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'-1.68%"
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).Value = "'0.8948"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'-3.35%"
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).Value = "'0.1644"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'-5.84%"
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,4).Value = "'0.07541"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'-5.08%"
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,4).Value = "'0.07962"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'-8.26%"
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,4).Value = "'0.02987"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'-3.82%"
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).Value = "'0.09994"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'-0.32%"
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,4).Value = "'0.001504"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'-1.34%"
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).Value = "'0.005858"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'-0.12%"
$ws.Cells.Item(16,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,4).Value = "'3.460"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'-0.02%"
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).Value = "'2.107"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'-6.60%"
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,5).Value = "'-0.29%"
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).Value = "'0.1298"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'-0.78%"
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).Value = "'4.269"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'-0.95%"
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,4).Value = "'0.2002"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'11.69%"
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).Value = "'0.04479"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'-2.74%"
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,4).Value = "'0.001210"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'-2.17%"
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).Value = "'0.004659"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'4.89%"
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,5).Value = "'-0.07%"
$ws.Cells.Item(27,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).Value = "'0.01636"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "'-4.56%"
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).Value = "'0.04350"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'-8.98%"
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).Value = "'0.007405"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'-1.27%"
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,5).Value = "'-3.82%"
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,5).Value = "'-15.23%"
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,4).Value = "'0.01024"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'-8.96%"
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,4).Value = "'0.00005846"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'-2.51%"
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,4).Value = "'0.00000000751"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'0.01%"
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,4).Value = "'2.216"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'170.05%"
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,5).Value = "'-11.46%"
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,5).Value = "'0.01%"
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,5).Value = "'0.01%"
$ws.Cells.Item(50,5).Style = "Normal"
